$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct table 2 data: Meningococcus row (row 4)
# Annual pre-vaccine deaths (C4): 0.00 -> 13.38
# Annual post-vaccine deaths (E4): 0.00 -> 2.60
# Use a leading apostrophe so the numeric-looking text is stored as text
# (matching the source workbook, where these figures are text strings),
# then reset the style back to Normal so no stray number-format/quote-
# prefix style is introduced.
$ws.Range("C4").Value = "'13.38"
$ws.Range("C4").Style = "Normal"

$ws.Range("E4").Value = "'2.60"
$ws.Range("E4").Style = "Normal"
